$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CPU part number updated: STM32F102CBU6 -> STM32F103CBU6
$ws.Range("D4").Value = "STM32F103CBU6"

# Row 7: 3.3V LDO replaced with a specific regulator part, with part number & DigiKey number added
$ws.Range("C7").Value = "IC REG LINEAR 3.3V 150MA SOT23-5"
$ws.Range("D7").Value = "LD3985M33R"
$ws.Range("E7").Value = "497-3504-1-ND"

# New row 15: 8MHz crystal (notes-only row, column G)
$ws.Range("G15").Value = "8MHz crystal"

# New row 16: USB Reenumeration transistor
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "TRANS NPN 40V 0.2A SOT23"
$ws.Range("D16").Value = "MMBT3904,215"
$ws.Range("E16").Value = "1727-4044-1-ND"
$ws.Range("F16").Value = "SOT23"
$ws.Range("G16").Value = "USB Reenumeration transistor"

# Column G widened to fit new longer text (36.1 maps to the stored XML width of 37)
$ws.Columns.Item(7).ColumnWidth = 36.1

# Update selection to match author's final cursor position
$ws.Range("C7").Select()
